# Fill in the "Score" column (3rd column) of the rubric table with the
# values that correspond to each row's "Points" entry. Each of these
# cells is currently an empty paragraph (just a centered jc pPr) and
# needs a single run containing the point value as text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$scores = @{
    2  = "1"
    3  = "1"
    6  = "2"
    7  = "2"
    8  = "2"
    9  = "2"
    10 = "2"
    11 = "2"
    12 = "2"
    13 = "2"
    14 = "2"
    15 = "2"
    16 = "2"
    18 = "3"
    19 = "3"
    25 = "30"
}

foreach ($rowIndex in $scores.Keys) {
    $cell = $t.Cell($rowIndex, 3)
    $cell.Range.Text = $scores[$rowIndex]
}
